# Regenerate merged AHB files:
#  - rename header labels *_old -> *_FV2310 and *_new -> *_FV2404
#  - wrap the used range in an Excel Table (ListObject)
#  - freeze the header row (top row)

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1) Rename the column headers in row 1.
for ($col = 1; $col -le 10; $col++) {
    $cell = $ws.Cells.Item(1, $col)
    $cell.Value2 = ($cell.Value2 -replace '_old$', '_FV2310')
}

for ($col = 12; $col -le 21; $col++) {
    $cell = $ws.Cells.Item(1, $col)
    $cell.Value2 = ($cell.Value2 -replace '_new$', '_FV2404')
}

# 2) Turn the used range into a proper Excel Table ("Table1").
$usedRange = $ws.Range("A1:U77")
$tbl = $ws.ListObjects.Add([Microsoft.Office.Interop.Excel.XlListObjectSourceType]::xlSrcRange, $usedRange, $null, [Microsoft.Office.Interop.Excel.XlYesNoGuess]::xlYes)
$tbl.Name = "Table1"

# 3) Freeze the top (header) row.
$ws.Range("A2").Select()
$excel.ActiveWindow.FreezePanes = $true
